$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-16 down to 10-17
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new weekly price record
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44580
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112003
$ws.Range("G9").Value = "Ajo"
$ws.Range("H9").Value = "Chino"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 18000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19000
$ws.Range("N9").Value = "$/caja 10 kilos"
$ws.Range("O9").Value = "China"
$ws.Range("P9").Value = 1900
$ws.Range("Q9").Value = 10
$ws.Range("R9").Value = "Hortaliza"
